$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $s = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $s
}

Set-TextValue $ws.Range("D2") "24.942.36"
Set-TextValue $ws.Range("E2") "  +2.30%  "
Set-TextValue $ws.Range("D3") "1.677.81"
Set-TextValue $ws.Range("E3") "  +1.75%  "
Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  -0.14%  "
Set-TextValue $ws.Range("D5") "328.38"
Set-TextValue $ws.Range("E5") "  +6.88%  "
Set-TextValue $ws.Range("D6") "0.9997"
Set-TextValue $ws.Range("E6") "  +0.07%  "
Set-TextValue $ws.Range("D7") "0.3659"
Set-TextValue $ws.Range("E7") "  +1.07%  "
Set-TextValue $ws.Range("D8") "46.88"
Set-TextValue $ws.Range("E8") "  -1.31%  "
Set-TextValue $ws.Range("D9") "0.3252"
Set-TextValue $ws.Range("E9") "  -0.55%  "
Set-TextValue $ws.Range("D10") "1.144"
Set-TextValue $ws.Range("E10") "  +2.16%  "
Set-TextValue $ws.Range("D11") "0.07082"
Set-TextValue $ws.Range("E11") "  +2.14%  "
Set-TextValue $ws.Range("E12") "  +0.06%  "
Set-TextValue $ws.Range("D13") "6.091"
Set-TextValue $ws.Range("E13") "  +2.80%  "
Set-TextValue $ws.Range("D14") "19.63"
Set-TextValue $ws.Range("E14") "  +2.17%  "
Set-TextValue $ws.Range("D15") "1.678.93"
Set-TextValue $ws.Range("E15") "  +1.47%  "
Set-TextValue $ws.Range("D16") "6.641"
Set-TextValue $ws.Range("E16") "  +0.72%  "
Set-TextValue $ws.Range("D17") "0.00001048"
Set-TextValue $ws.Range("E17") "  +0.99%  "
Set-TextValue $ws.Range("D18") "0.06593"
Set-TextValue $ws.Range("E18") "  +1.35%  "
Set-TextValue $ws.Range("D19") "0.9998"
Set-TextValue $ws.Range("E19") "  +0.07%  "
Set-TextValue $ws.Range("D20") "78.98"
Set-TextValue $ws.Range("E20") "  +3.36%  "
Set-TextValue $ws.Range("D21") "15.91"
Set-TextValue $ws.Range("E21") "  +1.58%  "
Set-TextValue $ws.Range("D22") "5.928"
Set-TextValue $ws.Range("E22") "  +0.26%  "
Set-TextValue $ws.Range("D23") "12.86"
Set-TextValue $ws.Range("E23") "  +4.10%  "
Set-TextValue $ws.Range("D24") "24.952.92"
Set-TextValue $ws.Range("E24") "  +2.38%  "
Set-TextValue $ws.Range("D25") "2.445"
Set-TextValue $ws.Range("E25") "  +0.70%  "
Set-TextValue $ws.Range("D26") "2.409"
Set-TextValue $ws.Range("E26") "  +2.96%  "
Set-TextValue $ws.Range("D27") "148.16"
Set-TextValue $ws.Range("E27") "  +1.18%  "
Set-TextValue $ws.Range("D28") "18.78"
Set-TextValue $ws.Range("E28") "  +2.61%  "
Set-TextValue $ws.Range("D29") "1.862.51"
Set-TextValue $ws.Range("E29") "  +1.32%  "
Set-TextValue $ws.Range("D30") "125.83"
Set-TextValue $ws.Range("E30") "  +1.23%  "
Set-TextValue $ws.Range("D31") "1.187"
Set-TextValue $ws.Range("E31") "  +1.60%  "
Set-TextValue $ws.Range("D32") "4.075"
Set-TextValue $ws.Range("E32") "  +0.86%  "
Set-TextValue $ws.Range("D33") "5.782"
Set-TextValue $ws.Range("E33") "  +3.06%  "
Set-TextValue $ws.Range("D34") "0.08493"
Set-TextValue $ws.Range("E34") "  +2.02%  "
Set-TextValue $ws.Range("D35") "1.645"
Set-TextValue $ws.Range("E35") "  -1.38%  "
Set-TextValue $ws.Range("D36") "12.31"
Set-TextValue $ws.Range("E36") "  +0.51%  "
Set-TextValue $ws.Range("D37") "5.178"
Set-TextValue $ws.Range("E37") "  -0.55%  "
Set-TextValue $ws.Range("D38") "0.02251"
Set-TextValue $ws.Range("E38") "  +2.30%  "
Set-TextValue $ws.Range("D39") "1.233"
Set-TextValue $ws.Range("E39") "  +2.31%  "
Set-TextValue $ws.Range("D40") "0.06032"
Set-TextValue $ws.Range("E40") "  +0.22%  "
Set-TextValue $ws.Range("D41") "0.2096"
Set-TextValue $ws.Range("D42") "8.229"
Set-TextValue $ws.Range("E42") "  +0.79%  "
Set-TextValue $ws.Range("D43") "0.9987"
Set-TextValue $ws.Range("E43") "  -0.06%  "
Set-TextValue $ws.Range("D44") "0.5965"
Set-TextValue $ws.Range("E44") "  +2.57%  "
Set-TextValue $ws.Range("D45") "13.68"
Set-TextValue $ws.Range("E45") "  +8.55%  "
Set-TextValue $ws.Range("D46") "3.842"
Set-TextValue $ws.Range("E46") "  +3.17%  "
Set-TextValue $ws.Range("E47") "  +3.80%  "
Set-TextValue $ws.Range("D48") "125.56"
Set-TextValue $ws.Range("E48") "  +3.31%  "
Set-TextValue $ws.Range("D49") "1.970"
Set-TextValue $ws.Range("E49") "  +1.90%  "
Set-TextValue $ws.Range("D50") "0.07026"
Set-TextValue $ws.Range("E50") "  +1.93%  "
Set-TextValue $ws.Range("D51") "1.191"
Set-TextValue $ws.Range("E51") "  +3.34%  "
